# Edit Steuertarife.xlsx: change H86:H113 increment from 100 to 1,
# and update the active sheet view (scroll position & selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Staatssteuer")

# Update "Einkommensinkrement" (column H) for rows 86-113 from 100 to 1
for ($r = 86; $r -le 113; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}

# Update the sheet view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 95
$ws.Range("N90").Select()
